# Invitados_con_QR.xlsx edit:
#   - The sheet used to carry two parallel "code / QR path" blocks:
#     columns C/D/E held an unused header-only trio ("Código Único",
#     "Enlace QR", "Asistencia") while the real per-invitee data lived in
#     columns F/G ("Código único" / "Ruta del QR"), with a duplicate
#     "Asistencia" header stranded in H1.
#   - The edit consolidates everything back down to a single 5-column
#     table: Nombre, Correo, Código Único, Ruta del QR, Asistencia -
#     reusing the original C1 header and moving the F/G data left into
#     C/D, and renaming D1 to "Ruta del QR" (taken from the old G1
#     header). The now-empty F:H columns are cleared out entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-point the D header at the text that used to live on G1 ("Ruta del
#    QR"); the old D1 text ("Enlace QR") is retired.
$ws.Range("D1").Value2 = "Ruta del QR"

# 2) Pull the real data (Código único + Ruta del QR) from F2:G98 up into
#    C2:D98, which were just empty header-only columns before.
$dataRange = $ws.Range("F2:G98")
$ws.Range("C2:D98").Value2 = $dataRange.Value2

# 3) The old F/G/H columns (stray header "Código único", the
#    already-copied "Ruta del QR" data, and the duplicate "Asistencia"
#    header in H1) are no longer needed - clear them out completely.
$ws.Range("F1:H98").Clear()

# 4) Widen the new D column (now holding the QR file paths) to roughly
#    the same width the data had back when it was column G.
$ws.Columns.Item(4).ColumnWidth = 23.3

# 5) Match the saved selection state.
$ws.Range("I5").Select()
